$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 760.375
$ws.Range("I6").Value = 760.375
$ws.Range("K6").Value = 2281.125
$ws.Range("M6").Value = -2169.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 57076.332
$ws.Range("I28").Value = 84349.086
$ws.Range("J28").Value = 2530.8333
$ws.Range("K28").Value = 84349.086
$ws.Range("L28").Value = 2530.8333
$ws.Range("M28").Value = -83864.086
$ws.Range("N28").Value = -3500.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 211.71428
$ws.Range("I33").Value = 211.71428
$ws.Range("K33").Value = 211.71428
$ws.Range("M33").Value = 17.28572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5316.5454
$ws.Range("I40").Value = 3569.9285
$ws.Range("J40").Value = 8373.125
$ws.Range("K40").Value = 3569.9285
$ws.Range("L40").Value = 8373.125
$ws.Range("M40").Value = -3394.9285
$ws.Range("N40").Value = -8723.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4219.231
$ws.Range("I98").Value = 1641.3334
$ws.Range("J98").Value = 4992.6
$ws.Range("K98").Value = 1641.3334
$ws.Range("L98").Value = 4992.6
$ws.Range("M98").Value = -143.3334
$ws.Range("N98").Value = -7988.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4219.231
$ws.Range("I122").Value = 1641.3334
$ws.Range("J122").Value = 4992.6
$ws.Range("K122").Value = 4924.0002
$ws.Range("L122").Value = 14977.8
$ws.Range("M122").Value = -2474.0002
$ws.Range("N122").Value = -19877.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2016.7941
$ws.Range("I132").Value = 1227.6207
$ws.Range("J132").Value = 6594
$ws.Range("K132").Value = 3682.8621
$ws.Range("L132").Value = 19782
$ws.Range("M132").Value = -1152.8621
$ws.Range("N132").Value = -24842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3575.0908
$ws.Range("I137").Value = 2157.853
$ws.Range("K137").Value = 6473.559
$ws.Range("M137").Value = -3923.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 5000
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4533

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 48888.5
$ws.Range("I39").Value = 50000
$ws.Range("J39").Value = 47777
$ws.Range("K39").Value = 50000
$ws.Range("L39").Value = 47777
$ws.Range("M39").Value = -49480
$ws.Range("N39").Value = -48817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6719.769
$ws.Range("I45").Value = 4085.8
$ws.Range("J45").Value = 15499.667
$ws.Range("K45").Value = 4085.8
$ws.Range("L45").Value = 15499.667
$ws.Range("M45").Value = -3708.8
$ws.Range("N45").Value = -16253.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1412.2142
$ws.Range("I97").Value = 1709.8
$ws.Range("J97").Value = 668.25
$ws.Range("K97").Value = 1709.8
$ws.Range("L97").Value = 668.25
$ws.Range("M97").Value = -1213.8
$ws.Range("N97").Value = -1660.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 40500
$ws.Range("J104").Value = 40500
$ws.Range("L104").Value = 40500
$ws.Range("N104").Value = -47488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 834606.2
$ws.Range("I110").Value = 834606.2
$ws.Range("K110").Value = 834606.2
$ws.Range("M110").Value = -832561.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10714.333
$ws.Range("I132").Value = 5357.25
$ws.Range("K132").Value = 16071.75
$ws.Range("M132").Value = -13541.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 50002230
$ws.Range("I20").Value = 76924680
$ws.Range("J20").Value = 3399.7144
$ws.Range("K20").Value = 76924680
$ws.Range("L20").Value = 3399.7144
$ws.Range("M20").Value = -76924433
$ws.Range("N20").Value = -3893.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4742
$ws.Range("I86").Value = 4717.476
$ws.Range("K86").Value = 4717.476
$ws.Range("M86").Value = -3594.476

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4742
$ws.Range("I89").Value = 4717.476
$ws.Range("K89").Value = 23587.38
$ws.Range("M89").Value = -17971.38

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3195.6
$ws.Range("I134").Value = 1961
$ws.Range("K134").Value = 5883
$ws.Range("M134").Value = -3348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3933.353
$ws.Range("I31").Value = 2219.7273
$ws.Range("K31").Value = 2219.7273
$ws.Range("M31").Value = -1924.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3933.353
$ws.Range("I34").Value = 2219.7273
$ws.Range("K34").Value = 2219.7273
$ws.Range("M34").Value = -2017.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1499.6
$ws.Range("I105").Value = 1463.8572
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1463.8572
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 283.1428000000001
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18310116
$ws.Range("I4").Value = 1118714.2
$ws.Range("K4").Value = 3356142.6
$ws.Range("M4").Value = -3356030.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 75884.60000000001
$ws.Range("J34").Value = 142225
$ws.Range("L34").Value = 426675
$ws.Range("N34").Value = -426843

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10842
$ws.Range("J39").Value = 15730
$ws.Range("L39").Value = 47190
$ws.Range("N39").Value = -47778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 30840
$ws.Range("I55").Value = 2600
$ws.Range("J55").Value = 49666.668
$ws.Range("K55").Value = 7800
$ws.Range("L55").Value = 149000.004
$ws.Range("M55").Value = -7623
$ws.Range("N55").Value = -149354.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9666.333000000001
$ws.Range("I102").Value = 9666.333000000001
$ws.Range("K102").Value = 9666.333000000001
$ws.Range("M102").Value = -8044.333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5232.38
$ws.Range("I122").Value = 5019.4375
$ws.Range("K122").Value = 15058.3125
$ws.Range("M122").Value = -12608.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1113488.5
$ws.Range("I132").Value = 1430828.1
$ws.Range("K132").Value = 4292484.300000001
$ws.Range("M132").Value = -4289954.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 458040.78
$ws.Range("I40").Value = 557938.75
$ws.Range("K40").Value = 557938.75
$ws.Range("M40").Value = -557802.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 710.0625
$ws.Range("I55").Value = 462.5
$ws.Range("J55").Value = 957.625
$ws.Range("K55").Value = 462.5
$ws.Range("L55").Value = 957.625
$ws.Range("M55").Value = -289.5
$ws.Range("N55").Value = -1303.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6533.6113
$ws.Range("I61").Value = 4461.4
$ws.Range("J61").Value = 9123.875
$ws.Range("K61").Value = 4461.4
$ws.Range("L61").Value = 9123.875
$ws.Range("M61").Value = -4259.4
$ws.Range("N61").Value = -9527.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7468.8184
$ws.Range("I68").Value = 7269.25
$ws.Range("J68").Value = 8001
$ws.Range("K68").Value = 7269.25
$ws.Range("L68").Value = 8001
$ws.Range("M68").Value = -6520.25
$ws.Range("N68").Value = -9499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7468.8184
$ws.Range("I71").Value = 7269.25
$ws.Range("J71").Value = 8001
$ws.Range("K71").Value = 36346.25
$ws.Range("L71").Value = 40005
$ws.Range("M71").Value = -32602.25
$ws.Range("N71").Value = -47493

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6533.6113
$ws.Range("I113").Value = 4461.4
$ws.Range("J113").Value = 9123.875
$ws.Range("K113").Value = 4461.4
$ws.Range("L113").Value = 9123.875
$ws.Range("M113").Value = -2291.4
$ws.Range("N113").Value = -13463.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3514.05
$ws.Range("I132").Value = 2361.5
$ws.Range("K132").Value = 7084.5
$ws.Range("M132").Value = -4554.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50584

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1424.5
$ws.Range("I100").Value = 999
$ws.Range("K100").Value = 1998
$ws.Range("M100").Value = -1457

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 49979.5
$ws.Range("J110").Value = 49979.5
$ws.Range("L110").Value = 49979.5
$ws.Range("N110").Value = -58159.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28575056
$ws.Range("I122").Value = 43480610
$ws.Range("K122").Value = 130441830
$ws.Range("M122").Value = -130439380
